$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Balance" running-total column (F) to the Semi Additive Formulae table
$ws.Range("F2").Value = "Balance"

# First data row just carries the amount across
$ws.Range("F3").Formula = "=E3"

# Second row adds current amount to previous balance (not yet part of the shared formula)
$ws.Range("F4").Formula = "=E4+F3"

# Remaining rows share one formula definition (F5:F41)
$ws.Range("F5:F41").Formula = "=E5+F4"

# Restore the selection to the cell that was active when the workbook was saved
$ws.Range("F7").Select()
